$wb = $excel.ActiveWorkbook

# Row 33 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 27027878
$ws.Range("I33").Value = 621.0303
$ws.Range("J33").Value = 250002750
$ws.Range("K33").Value = 621.0303
$ws.Range("L33").Value = 250002750
$ws.Range("M33").Value = -392.0303
$ws.Range("N33").Value = -250003208

# Row 76 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3083.2856
$ws.Range("I76").Value = 3082.679
$ws.Range("J76").Value = 3099.6667
$ws.Range("K76").Value = 3082.679
$ws.Range("L76").Value = 3099.6667
$ws.Range("M76").Value = -2767.679
$ws.Range("N76").Value = -3729.6667

# Row 79 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3083.2856
$ws.Range("I79").Value = 3082.679
$ws.Range("J79").Value = 3099.6667
$ws.Range("K79").Value = 3082.679
$ws.Range("L79").Value = 3099.6667
$ws.Range("M79").Value = -1990.679
$ws.Range("N79").Value = -5283.6667

# Row 100 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2826.3635
$ws.Range("I100").Value = 2190
$ws.Range("J100").Value = 3190
$ws.Range("K100").Value = 2190
$ws.Range("L100").Value = 3190
$ws.Range("M100").Value = -1649
$ws.Range("N100").Value = -4272

# Row 129 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 933
$ws.Range("I129").Value = 434.22223
$ws.Range("J129").Value = 1232.2667
$ws.Range("K129").Value = 1302.66669
$ws.Range("L129").Value = 3696.800099999999
$ws.Range("M129").Value = 3697.33331
$ws.Range("N129").Value = -13696.8001

# Row 64 (ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67 (ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 105 (BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2468.0667
$ws.Range("I105").Value = 2010
$ws.Range("J105").Value = 2500.7856
$ws.Range("K105").Value = 2010
$ws.Range("L105").Value = 2500.7856
$ws.Range("M105").Value = -263
$ws.Range("N105").Value = -5994.7856

# Row 107 (BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4238.8
$ws.Range("I107").Value = 3074.5
$ws.Range("J107").Value = 6567.4
$ws.Range("K107").Value = 3074.5
$ws.Range("L107").Value = 6567.4
$ws.Range("M107").Value = -1154.5
$ws.Range("N107").Value = -10407.4

# Row 31 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1316.8889
$ws.Range("I31").Value = 1316.8889
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1316.8889
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1021.8889
$ws.Range("N31").ClearContents()

# Row 34 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1316.8889
$ws.Range("I34").Value = 1316.8889
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1316.8889
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1114.8889
$ws.Range("N34").ClearContents()

# Row 52 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 30337.8
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 30337.8
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 30337.8
$ws.Range("N52").Value = -30925.8

# Row 3 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3117.3333
$ws.Range("I3").Value = 2366.1538
$ws.Range("J3").Value = 8000
$ws.Range("K3").Value = 7098.4614
$ws.Range("L3").Value = 24000
$ws.Range("M3").Value = -6986.4614
$ws.Range("N3").Value = -24224

# Row 133 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 7210
$ws.Range("I133").Value = 4420
$ws.Range("J133").Value = 10000
$ws.Range("K133").Value = 13260
$ws.Range("L133").Value = 30000
$ws.Range("M133").Value = -8200
$ws.Range("N133").Value = -40120

# Row 134 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 6763.3335
$ws.Range("I134").Value = 5190
$ws.Range("J134").Value = 7550
$ws.Range("K134").Value = 15570
$ws.Range("L134").Value = 22650
$ws.Range("M134").Value = -10500
$ws.Range("N134").Value = -32790

# Row 136 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 62504724
$ws.Range("I136").Value = 100001560
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 300004680
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -299999580
$ws.Range("N136").Value = -40200

# Row 137 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 52057976
$ws.Range("I137").Value = 6000
$ws.Range("J137").Value = 56395640
$ws.Range("K137").Value = 18000
$ws.Range("L137").Value = 169186920
$ws.Range("M137").Value = -12900
$ws.Range("N137").Value = -169197120

# Row 138 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2215.75
$ws.Range("I138").Value = 2276.6667
$ws.Range("J138").Value = 2033
$ws.Range("K138").Value = 6830.000100000001
$ws.Range("L138").Value = 6099
$ws.Range("M138").Value = -1690.000100000001
$ws.Range("N138").Value = -16379

# Row 139 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2406.6667
$ws.Range("I139").Value = 1728.5714
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 5185.7142
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -45.71420000000035
$ws.Range("N139").Value = -19280

# Row 140 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 4393.6665
$ws.Range("I140").Value = 1433.3334
$ws.Range("J140").Value = 5380.4443
$ws.Range("K140").Value = 4300.0002
$ws.Range("L140").Value = 16141.3329
$ws.Range("M140").Value = 879.9997999999996
$ws.Range("N140").Value = -26501.3329

# Row 70 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4853.283
$ws.Range("I70").Value = 4710.9614
$ws.Range("J70").Value = 5778.375
$ws.Range("K70").Value = 4710.9614
$ws.Range("L70").Value = 5778.375
$ws.Range("M70").Value = -4440.9614
$ws.Range("N70").Value = -6318.375

# Row 73 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4853.283
$ws.Range("I73").Value = 4710.9614
$ws.Range("J73").Value = 5778.375
$ws.Range("K73").Value = 4710.9614
$ws.Range("L73").Value = 5778.375
$ws.Range("M73").Value = -3774.9614
$ws.Range("N73").Value = -7650.375

# Row 113 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1857.5
$ws.Range("I113").Value = 1363.3
$ws.Range("J113").Value = 2269.3333
$ws.Range("K113").Value = 1363.3
$ws.Range("L113").Value = 2269.3333
$ws.Range("M113").Value = 806.7
$ws.Range("N113").Value = -6609.3333

# Row 11 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 31002.334
$ws.Range("I11").Value = 11000
$ws.Range("J11").Value = 41003.5
$ws.Range("K11").Value = 11000
$ws.Range("L11").Value = 41003.5
$ws.Range("M11").Value = -10860
$ws.Range("N11").Value = -41283.5

# Row 55 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 800.9167
$ws.Range("I55").Value = 222.5
$ws.Range("J55").Value = 1090.125
$ws.Range("K55").Value = 222.5
$ws.Range("L55").Value = 1090.125
$ws.Range("M55").Value = -49.5
$ws.Range("N55").Value = -1436.125

# Row 61 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1849.95
$ws.Range("I61").Value = 1421.3572
$ws.Range("J61").Value = 2850
$ws.Range("K61").Value = 1421.3572
$ws.Range("L61").Value = 2850
$ws.Range("M61").Value = -1219.3572
$ws.Range("N61").Value = -3254

# Row 113 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1849.95
$ws.Range("I113").Value = 1421.3572
$ws.Range("J113").Value = 2850
$ws.Range("K113").Value = 1421.3572
$ws.Range("L113").Value = 2850
$ws.Range("M113").Value = 748.6428000000001
$ws.Range("N113").Value = -7190

# Row 122 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 41952.633
$ws.Range("I122").Value = 1969.6666
$ws.Range("J122").Value = 101927.086
$ws.Range("K122").Value = 5908.9998
$ws.Range("L122").Value = 305781.258
$ws.Range("M122").Value = -3458.9998
$ws.Range("N122").Value = -310681.258

# Row 126 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1497.826
$ws.Range("I126").Value = 1030.5454
$ws.Range("J126").Value = 1926.1666
$ws.Range("K126").Value = 3091.6362
$ws.Range("L126").Value = 5778.4998
$ws.Range("M126").Value = -621.6361999999999
$ws.Range("N126").Value = -10718.4998

# Row 132 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3199.7026
$ws.Range("I132").Value = 3095.48
$ws.Range("J132").Value = 3416.8333
$ws.Range("K132").Value = 9286.440000000001
$ws.Range("L132").Value = 10250.4999
$ws.Range("M132").Value = -6756.440000000001
$ws.Range("N132").Value = -15310.4999

# Row 136 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3133.7856
$ws.Range("I136").Value = 3297.1904
$ws.Range("J136").Value = 2643.5715
$ws.Range("K136").Value = 9891.5712
$ws.Range("L136").Value = 7930.7145
$ws.Range("M136").Value = -7341.5712
$ws.Range("N136").Value = -13030.7145
